# Bump the ObjTables schema header metadata embedded in this "Address book"
# workbook:
#   objTablesVersion: 0.0.9 -> 1.0.0
#   date:             2020-04-28 15:09:22 -> 2020-05-29 00:17:37
#
# The header strings live in cell A1 (table-of-contents sheet also has a
# second header in A2) of every worksheet. The sheets are protected, so each
# one has to be unprotected before the edit and re-protected afterwards.

$wb = $excel.ActiveWorkbook

$oldVersion = "0.0.9"
$newVersion = "1.0.0"
$oldDate = "2020-04-28 15:09:22"
$newDate = "2020-05-29 00:17:37"

foreach ($ws in $wb.Worksheets) {
    $ws.Unprotect()

    foreach ($row in 1, 2) {
        $cell = $ws.Cells.Item($row, 1)
        $text = $cell.Text

        if ($text -like "!!ObjTables*" -or $text -like "!!!ObjTables*") {
            $updated = $text.Replace($oldVersion, $newVersion).Replace($oldDate, $newDate)
            if ($updated -ne $text) {
                $cell.Value = $updated
            }
        }
    }

    $ws.Protect()
}
